$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32 (shifts existing rows 32..151 down to 33..152,
# preserving their values/formatting, and extends the used range to A1:T152).
$ws.Rows(32).Insert()

# Populate the newly inserted row 32 with the new daily price record.
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C32").Value = "Ñuble"
$ws.Range("D32").Value = [DateTime]"2023-06-30"
$ws.Range("D32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100108
$ws.Range("H32").Value = "Tropicales y subtropicales"
$ws.Range("I32").Value = 100108002
$ws.Range("J32").Value = "Mango"
$ws.Range("K32").Value = "Sin especificar"
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 40
$ws.Range("N32").Value = 9000
$ws.Range("O32").Value = 9000
$ws.Range("P32").Value = 9000
$ws.Range("Q32").Value = '$/bandeja 4 kilos'
$ws.Range("R32").Value = "Perú"
$ws.Range("S32").Value = 2250
$ws.Range("T32").Value = 4
